# Removing the invalid test cases for watchlist
#
# The "Test Cases" sheet (sheet1 / "Test Cases") tracked watchlist test
# cases TestCase_E1..E10. The invalid ones (E5-E10) are removed, and the
# remaining four (E1-E4) get their Description/Results updated to reflect
# the new "watch/unwatch an Article" wording.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Drop the rows for the removed/invalid test cases (TestCase_E5..E10).
$ws.Rows("6:11").Delete() | Out-Null

# Update the remaining rows' Description (C) and Results (E) columns.
$ws.Range("C2").Value = "Verify that user is able to watch an Article from ALL content search results page"
$ws.Range("E2").Value = "SKIP"

$ws.Range("C3").Value = "Verify that user is able to watch an Article from Record View page"
$ws.Range("E3").Value = "PASS"

$ws.Range("C4").Value = "Verify that user is able to unwatch an Article from watchlist page"
$ws.Range("E4").Value = "PASS"

$ws.Range("C5").Value = "Verify that user is able to unwatch an Article from ALL content search results page"
$ws.Range("E5").Value = "SKIP"

# Re-fit column C now that the longest description text has changed.
$ws.Columns("C").ColumnWidth = 104.5

# Match the author's final selection/viewport on the shrunk sheet.
$ws.Range("D5").Select() | Out-Null
